$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: header address line
#   "AV: Abdo Jaud Feres, " -> "AV: Abdo Jauid Feres, "
#   (split into 3 runs: "AV: Abdo Jau" / "i" / "d Feres, ", same bold rPr)
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute("Jaud")
$start1 = $rng1.Start
$d.Range($start1 + 3, $start1 + 3).InsertAfter("i")
$iRun1 = $d.Range($start1 + 3, $start1 + 4)
$iRun1.Bold = 0
$iRun1.Bold = 1

# ---------------------------------------------------------------------------
# Change 2: same paragraph, mark "Angelo" as its own run
#   "165 - Jardim Angelo Passuelo " ->
#   "165 - Jardim " + "Angelo" + " Passuelo " (3 runs, identical rPr)
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("Angelo")
$angeloRun = $d.Range($rng2.Start, $rng2.End)
$angeloRun.Bold = 0
$angeloRun.Bold = 1

# ---------------------------------------------------------------------------
# Change 3: body paragraph
#   "...Av. Abdo Jaud Feres, 165..." -> "...Av. Abdo Jauid Feres, 165..."
#   (split into 3 runs: ".. Jau" / "i" / "d Feres, 165.."; this run has no
#    explicit Bold/Italic in rPr, so round-trip through wdUndefined to avoid
#    leaving a stray attribute behind)
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("Jaud")
$start3 = $rng3.Start
$d.Range($start3 + 3, $start3 + 3).InsertAfter("i")
$iRun3 = $d.Range($start3 + 3, $start3 + 4)
$iRun3.Bold = 1
$iRun3.Bold = 9999999

Write-Host "Done"
